$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC (47 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 974.8333
$ws.Range("I5").Value = 16.666666
$ws.Range("J5").Value = 1933
$ws.Range("K5").Value = 16.666666
$ws.Range("L5").Value = 1933
$ws.Range("M5").Value = 98.33333400000001
$ws.Range("N5").Value = -2163
$ws.Range("H41").Value = 360.94116
$ws.Range("I41").Value = 88.2
$ws.Range("J41").Value = 474.58334
$ws.Range("K41").Value = 88.2
$ws.Range("L41").Value = 474.58334
$ws.Range("M41").Value = 351.8
$ws.Range("N41").Value = -1354.58334
$ws.Range("H76").Value = 2927239.8
$ws.Range("I76").Value = 3461.5386
$ws.Range("K76").Value = 3461.5386
$ws.Range("M76").Value = -3146.5386
$ws.Range("H79").Value = 2927239.8
$ws.Range("I79").Value = 3461.5386
$ws.Range("K79").Value = 3461.5386
$ws.Range("M79").Value = -2369.5386
$ws.Range("H129").Value = 193301.53
$ws.Range("I129").Value = 347.5
$ws.Range("J129").Value = 209381.05
$ws.Range("K129").Value = 1042.5
$ws.Range("L129").Value = 628143.1499999999
$ws.Range("M129").Value = 3957.5
$ws.Range("N129").Value = -638143.1499999999
$ws.Range("H132").Value = 4107.4346
$ws.Range("I132").Value = 4498.55
$ws.Range("K132").Value = 13495.65
$ws.Range("M132").Value = -10965.65
$ws.Range("H135").Value = 19232204
$ws.Range("I135").Value = 1104.1333
$ws.Range("J135").Value = 45456430
$ws.Range("K135").Value = 9937.199699999999
$ws.Range("L135").Value = 409107870
$ws.Range("M135").Value = -7402.199699999999
$ws.Range("N135").Value = -409112940
$ws.Range("H137").Value = 1135.9117
$ws.Range("I137").Value = 1066.3334
$ws.Range("J137").Value = 1302.9
$ws.Range("K137").Value = 3199.0002
$ws.Range("L137").Value = 3908.7
$ws.Range("M137").Value = -649.0001999999999
$ws.Range("N137").Value = -9008.700000000001

# ---- Worksheet: ARM (49 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 682.8333
$ws.Range("I2").Value = 665.62964
$ws.Range("K2").Value = 665.62964
$ws.Range("M2").Value = -552.62964
$ws.Range("H32").Value = 6580.6313
$ws.Range("I32").Value = 5361.547
$ws.Range("J32").Value = 13082.417
$ws.Range("K32").Value = 5361.547
$ws.Range("L32").Value = 13082.417
$ws.Range("M32").Value = -5074.547
$ws.Range("N32").Value = -13656.417
$ws.Range("H58").Value = 22000
$ws.Range("J58").Value = 22000
$ws.Range("L58").Value = 22000
$ws.Range("N58").Value = -22860
$ws.Range("H63").Value = 3908143.8
$ws.Range("I63").Value = 2164.2856
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 2164.2856
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -1478.2856
$ws.Range("N63").Value = -31251372
$ws.Range("H66").Value = 3908143.8
$ws.Range("I66").Value = 2164.2856
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 10821.428
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -7389.428
$ws.Range("N66").Value = -156256864
$ws.Range("H116").Value = 682.8333
$ws.Range("I116").Value = 665.62964
$ws.Range("K116").Value = 665.62964
$ws.Range("M116").Value = 1628.37036
$ws.Range("H121").Value = 27375
$ws.Range("J121").Value = 27375
$ws.Range("L121").Value = 27375
$ws.Range("N121").Value = -30869
$ws.Range("H122").Value = 1278.96
$ws.Range("I122").Value = 1298.9166
$ws.Range("K122").Value = 3896.7498
$ws.Range("M122").Value = -1446.7498
$ws.Range("H132").Value = 12290.25
$ws.Range("I132").Value = 1662.0513
$ws.Range("K132").Value = 4986.1539
$ws.Range("M132").Value = -2456.1539
$ws.Range("H133").Value = 52603
$ws.Range("J133").Value = 52603
$ws.Range("L133").Value = 52603
$ws.Range("N133").Value = -57663

# ---- Worksheet: BSM (29 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 682.8333
$ws.Range("I3").Value = 665.62964
$ws.Range("K3").Value = 665.62964
$ws.Range("M3").Value = -551.62964
$ws.Range("H64").Value = 271.64706
$ws.Range("I64").Value = 143.71428
$ws.Range("J64").Value = 361.2
$ws.Range("K64").Value = 143.71428
$ws.Range("L64").Value = 361.2
$ws.Range("M64").Value = 81.28572
$ws.Range("N64").Value = -811.2
$ws.Range("H67").Value = 271.64706
$ws.Range("I67").Value = 143.71428
$ws.Range("J67").Value = 361.2
$ws.Range("K67").Value = 143.71428
$ws.Range("L67").Value = 361.2
$ws.Range("M67").Value = 636.28572
$ws.Range("N67").Value = -1921.2
$ws.Range("H99").Value = 1970.3334
$ws.Range("I99").Value = 1950
$ws.Range("K99").Value = 1950
$ws.Range("M99").Value = -452
$ws.Range("H134").Value = 3217.7346
$ws.Range("I134").Value = 3240.475
$ws.Range("J134").Value = 3116.6667
$ws.Range("K134").Value = 9721.424999999999
$ws.Range("L134").Value = 9350.000100000001
$ws.Range("M134").Value = -7186.424999999999
$ws.Range("N134").Value = -14420.0001

# ---- Worksheet: CRP (34 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 20000
$ws.Range("J41").Value = 20000
$ws.Range("L41").Value = 20000
$ws.Range("N41").Value = -20856
$ws.Range("H51").Value = 34500
$ws.Range("J51").Value = 34500
$ws.Range("L51").Value = 34500
$ws.Range("N51").Value = -35972
$ws.Range("H61").Value = 34500
$ws.Range("J61").Value = 34500
$ws.Range("L61").Value = 34500
$ws.Range("N61").Value = -35196
$ws.Range("H107").Value = 1084.8276
$ws.Range("I107").Value = 396.05264
$ws.Range("J107").Value = 2393.5
$ws.Range("K107").Value = 396.05264
$ws.Range("L107").Value = 2393.5
$ws.Range("M107").Value = 1523.94736
$ws.Range("N107").Value = -6233.5
$ws.Range("H122").Value = 1631.8182
$ws.Range("I122").Value = 1493.75
$ws.Range("K122").Value = 4481.25
$ws.Range("M122").Value = -2031.25
$ws.Range("H132").Value = 2381.4443
$ws.Range("I132").Value = 1778.7241
$ws.Range("J132").Value = 4878.4287
$ws.Range("K132").Value = 5336.1723
$ws.Range("L132").Value = 14635.2861
$ws.Range("M132").Value = -2806.1723
$ws.Range("N132").Value = -19695.2861
$ws.Range("H134").Value = 1066.8334
$ws.Range("I134").Value = 880.2
$ws.Range("K134").Value = 2640.6
$ws.Range("M134").Value = -105.6000000000004

# ---- Worksheet: CUL (4 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 684.85
$ws.Range("J131").Value = 711.9231
$ws.Range("L131").Value = 2135.7693
$ws.Range("N131").Value = -12215.7693

# ---- Worksheet: GSM (15 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 21001.428
$ws.Range("I46").Value = 20000
$ws.Range("J46").Value = 21168.334
$ws.Range("K46").Value = 20000
$ws.Range("L46").Value = 21168.334
$ws.Range("M46").Value = -19844
$ws.Range("N46").Value = -21480.334
$ws.Range("H131").Value = 46326
$ws.Range("J131").Value = 46326
$ws.Range("L131").Value = 46326
$ws.Range("N131").Value = -56406
$ws.Range("H135").Value = 39770
$ws.Range("J135").Value = 39770
$ws.Range("L135").Value = 39770
$ws.Range("N135").Value = -49910

# ---- Worksheet: LTW (21 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2566.8125
$ws.Range("I22").Value = 3333.25
$ws.Range("J22").Value = 267.5
$ws.Range("K22").Value = 3333.25
$ws.Range("L22").Value = 267.5
$ws.Range("M22").Value = -3038.25
$ws.Range("N22").Value = -857.5
$ws.Range("H27").Value = 2566.8125
$ws.Range("I27").Value = 3333.25
$ws.Range("J27").Value = 267.5
$ws.Range("K27").Value = 3333.25
$ws.Range("L27").Value = 267.5
$ws.Range("M27").Value = -3226.25
$ws.Range("N27").Value = -481.5
$ws.Range("H46").Value = 3022.2222
$ws.Range("I46").Value = 3716.6667
$ws.Range("J46").Value = 1633.3334
$ws.Range("K46").Value = 3716.6667
$ws.Range("L46").Value = 1633.3334
$ws.Range("M46").Value = -3528.6667
$ws.Range("N46").Value = -2009.3334

# ---- Worksheet: WVR (8 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 35000
$ws.Range("J119").Value = 35000
$ws.Range("L119").Value = 35000
$ws.Range("N119").Value = -44676
$ws.Range("H126").Value = 1839.48
$ws.Range("I126").Value = 1374.35
$ws.Range("K126").Value = 4123.049999999999
$ws.Range("M126").Value = -1653.049999999999
